# Apply updated "Actual Price" (J3) and purchase-lot (B/D columns) values
# for each crypto-asset worksheet, per the commit diff.
$wb = $excel.ActiveWorkbook

# Sheet 1: ETH
$ws = $wb.Worksheets.Item("ETH")
$ws.Range("J3").Value = [double]"2779.153168932849"
$ws.Range("B12").Value = [double]"0.00724113"
$ws.Range("B36").Value = [double]"0.02510273"
$ws.Range("D36").Value = [double]"44.1"
$ws.Range("B40").Value = [double]"0.05693671"
$ws.Range("D40").Value = [double]"107.05"

# Sheet 2: BTC
$ws = $wb.Worksheets.Item("BTC")
$ws.Range("J3").Value = [double]"51803.15604324524"
$ws.Range("B6").Value = [double]"0.00035561"
$ws.Range("B24").Value = [double]"0.00164913"
$ws.Range("D24").Value = [double]"44.1"
$ws.Range("B34").Value = [double]"0.00206961"
$ws.Range("D34").Value = [double]"60.95"

# Sheet 4: POLIS
$ws = $wb.Worksheets.Item("POLIS")
$ws.Range("J3").Value = [double]"0.4429032363911423"

# Sheet 5: ATLAS
$ws = $wb.Worksheets.Item("ATLAS")
$ws.Range("J3").Value = [double]"0.004918668355184971"

# Sheet 7: ACE
$ws = $wb.Worksheets.Item("ACE")
$ws.Range("J3").Value = [double]"10.32832620568351"
$ws.Range("B6").Value = [double]"2.48e-05"

# Sheet 8: ADA
$ws = $wb.Worksheets.Item("ADA")
$ws.Range("J3").Value = [double]"0.5926761381731338"
$ws.Range("B6").Value = [double]"0.7804502"
$ws.Range("B7").Value = [double]"122.48898156"
$ws.Range("D7").Value = [double]"44.1"

# Sheet 9: ALGO
$ws = $wb.Worksheets.Item("ALGO")
$ws.Range("J3").Value = [double]"0.1917923950985028"
$ws.Range("B6").Value = [double]"0.58399617"

# Sheet 10: APE
$ws = $wb.Worksheets.Item("APE")
$ws.Range("J3").Value = [double]"1.634783704549915"
$ws.Range("B5").Value = [double]"16.10617001"
$ws.Range("D5").Value = [double]"44.1"
$ws.Range("B6").Value = [double]"0.59985942"

# Sheet 11: ATOM
$ws = $wb.Worksheets.Item("ATOM")
$ws.Range("J3").Value = [double]"10.12966350442941"
$ws.Range("B7").Value = [double]"0.02906058"

# Sheet 12: AVAX
$ws = $wb.Worksheets.Item("AVAX")
$ws.Range("J3").Value = [double]"39.9377902126713"
$ws.Range("B5").Value = [double]"2.64960441"
$ws.Range("D5").Value = [double]"44.1"
$ws.Range("B6").Value = [double]"0.0167087"

# Sheet 13: AMP
$ws = $wb.Worksheets.Item("AMP")
$ws.Range("J3").Value = [double]"0.003703883154469359"

# Sheet 14: BNB
$ws = $wb.Worksheets.Item("BNB")
$ws.Range("J3").Value = [double]"355.7644748899322"
$ws.Range("B10").Value = [double]"0.00271795"
$ws.Range("B12").Value = [double]"0.15384525"
$ws.Range("D12").Value = [double]"44.1"

# Sheet 15: DOGE
$ws = $wb.Worksheets.Item("DOGE")
$ws.Range("J3").Value = [double]"0.08515996132831538"
$ws.Range("B6").Value = [double]"0.29102607"

# Sheet 16: DOT
$ws = $wb.Worksheets.Item("DOT")
$ws.Range("J3").Value = [double]"7.600877611838182"
$ws.Range("B5").Value = [double]"7.709423"
$ws.Range("D5").Value = [double]"44.1"
$ws.Range("B6").Value = [double]"0.07949792999999999"

# Sheet 17: EGLD
$ws = $wb.Worksheets.Item("EGLD")
$ws.Range("J3").Value = [double]"58.69927009932972"
$ws.Range("B6").Value = [double]"0.00299596"

# Sheet 18: GRT
$ws = $wb.Worksheets.Item("GRT")
$ws.Range("J3").Value = [double]"0.1918749210975464"

# Sheet 19: ICP
$ws = $wb.Worksheets.Item("ICP")
$ws.Range("J3").Value = [double]"12.99626879109305"
$ws.Range("B6").Value = [double]"0.00235763"

# Sheet 20: KAVA
$ws = $wb.Worksheets.Item("KAVA")
$ws.Range("J3").Value = [double]"0.7463328340742028"

# Sheet 21: LDO
$ws = $wb.Worksheets.Item("LDO")
$ws.Range("J3").Value = [double]"3.219422606909083"
$ws.Range("B6").Value = [double]"0.02016419"

# Sheet 22: LINK
$ws = $wb.Worksheets.Item("LINK")
$ws.Range("J3").Value = [double]"19.39722478869377"
$ws.Range("B6").Value = [double]"0.00247963"

# Sheet 23: LTC
$ws = $wb.Worksheets.Item("LTC")
$ws.Range("J3").Value = [double]"69.92622946464374"
$ws.Range("B6").Value = [double]"0.00133469"

# Sheet 24: LUNA
$ws = $wb.Worksheets.Item("LUNA")
$ws.Range("J3").Value = [double]"0.7078832392471571"
$ws.Range("B6").Value = [double]"0.05838977"

# Sheet 25: LUNC
$ws = $wb.Worksheets.Item("LUNC")
$ws.Range("J3").Value = [double]"0.0001250348425662677"
$ws.Range("B18").Value = [double]"5020.47704474"

# Sheet 26: MATIC
$ws = $wb.Worksheets.Item("MATIC")
$ws.Range("J3").Value = [double]"0.9349332125424628"
$ws.Range("B6").Value = [double]"0.32893872"
$ws.Range("B7").Value = [double]"48.53469011"
$ws.Range("D7").Value = [double]"44.1"

# Sheet 27: MEME
$ws = $wb.Worksheets.Item("MEME")
$ws.Range("J3").Value = [double]"0.0260155208868841"
$ws.Range("B6").Value = [double]"0.06762752"

# Sheet 28: MINA
$ws = $wb.Worksheets.Item("MINA")
$ws.Range("J3").Value = [double]"1.372767205897603"
$ws.Range("B6").Value = [double]"0.3498678"

# Sheet 29: NEAR
$ws = $wb.Worksheets.Item("NEAR")
$ws.Range("J3").Value = [double]"3.232064450321663"
$ws.Range("B6").Value = [double]"23.8991302"
$ws.Range("D6").Value = [double]"44.1"
$ws.Range("B7").Value = [double]"0.10290843"

# Sheet 30: SEI
$ws = $wb.Worksheets.Item("SEI")
$ws.Range("J3").Value = [double]"0.9307697269216131"
$ws.Range("B6").Value = [double]"0.07634290000000001"

# Sheet 31: SHIB
$ws = $wb.Worksheets.Item("SHIB")
$ws.Range("J3").Value = [double]"9.751985603361463e-06"
$ws.Range("B6").Value = [double]"275.37"

# Sheet 32: SHPING
$ws = $wb.Worksheets.Item("SHPING")
$ws.Range("J3").Value = [double]"0.004853732751778108"

# Sheet 33: SOL
$ws = $wb.Worksheets.Item("SOL")
$ws.Range("J3").Value = [double]"109.2429020961454"
$ws.Range("B17").Value = [double]"0.06468936"
$ws.Range("B18").Value = [double]"1.91528865"
$ws.Range("D18").Value = [double]"44.1"

# Sheet 34: TRX
$ws = $wb.Worksheets.Item("TRX")
$ws.Range("J3").Value = [double]"0.1320521835004282"
$ws.Range("B6").Value = [double]"0.26491802"

# Sheet 35: UNI
$ws = $wb.Worksheets.Item("UNI")
$ws.Range("J3").Value = [double]"7.417267145322007"
$ws.Range("B6").Value = [double]"0.00274524"

# Sheet 36: XRP
$ws = $wb.Worksheets.Item("XRP")
$ws.Range("J3").Value = [double]"0.5627741005570429"
$ws.Range("B6").Value = [double]"0.86409945"

# Sheet 37: TIA
$ws = $wb.Worksheets.Item("TIA")
$ws.Range("J3").Value = [double]"18.02564972583052"
$ws.Range("B6").Value = [double]"0.00399676"

# Sheet 38: DYDX
$ws = $wb.Worksheets.Item("DYDX")
$ws.Range("J3").Value = [double]"3.094561909438139"
$ws.Range("B6").Value = [double]"0.00079147"
